# Update average_county_temperature values (column AD) with NOAA data
# on Sheet1, per the commit "Updated temperature with NOAA data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-13 -> 21.28240740740739
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 30).Value = 21.28240740740739
}

# Rows 38-41 -> 19.65277777777778
for ($r = 38; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 19.65277777777778
}

# Rows 42-53 -> 13.75752314814816
for ($r = 42; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 13.75752314814816
}

# Rows 70-73 -> 19.65277777777778
for ($r = 70; $r -le 73; $r++) {
    $ws.Cells.Item($r, 30).Value = 19.65277777777778
}
